# Add setup of LIS2DH and logging without persistence to energy consumption
# New rows 20-29 on Tabelle1: "setup" (A20:A24) and "logging" (A25:A29) actions,
# both with B="without" (logging without persistence), following the existing
# Aktion / Modus / Resolution / Samplingrate / Spannung / Average Current /
# Dauer / Energie layout and the E*F*G/1000 energy formula used by the rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$rows = @(
    @{ Row=20; A="setup";   B="without"; C=10; D=10;  E=600; F=3; G=29 },
    @{ Row=21; A="setup";   B="without"; C=10; D=25;  E=600; F=3; G=31 },
    @{ Row=22; A="setup";   B="without"; C=10; D=50;  E=600; F=3; G=37 },
    @{ Row=23; A="setup";   B="without"; C=10; D=100; E=600; F=3; G=48 },
    @{ Row=24; A="setup";   B="without"; C=10; D=200; E=600; F=3; G=70 },
    @{ Row=25; A="logging"; B="without"; C=10; D=10;  E=600; F=3; G=31 },
    @{ Row=26; A="logging"; B="without"; C=10; D=25;  E=600; F=3; G=35 },
    @{ Row=27; A="logging"; B="without"; C=10; D=50;  E=600; F=3; G=42 },
    @{ Row=28; A="logging"; B="without"; C=10; D=100; E=600; F=3; G=59 },
    @{ Row=29; A="logging"; B="without"; C=10; D=200; E=600; F=3; G=91 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Formula = "=E" + $n + "*F" + $n + "*G" + $n + "/1000"
}

# Move / show the selection where the user left off after entering the new rows
$ws.Range("A30").Select() | Out-Null
